$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update header row (row 1) values for columns B:E
$ws.Range("B1").Value = 16
$ws.Range("C1").Value = 20
$ws.Range("D1").Value = 16
$ws.Range("E1").Value = 20

# Update row 2 values for columns B:E
$ws.Range("B2").Value = 5.2405407220496185
$ws.Range("C2").Value = 7.0105172500430539
$ws.Range("D2").Value = 11.218901824499449
$ws.Range("E2").Value = 10.457440156659771

# Update row 3 values for columns B:E
$ws.Range("B3").Value = 4.6333683627870021
$ws.Range("C3").Value = 7.34948576015063
$ws.Range("D3").Value = 7.02036264270017
$ws.Range("E3").Value = 11.396259373873145

# Shrink the selected range to match the updated data extent
$ws.Range("B1:E3").Select()
